# Auto-generated update script for khl_referees_stats workbook
# Updates per-referee season stats and refreshes the AA "last updated" timestamp
$wb = $excel.ActiveWorkbook

$timestamp = "2025-11-21 03:03:53"

$wsMain = $wb.Worksheets.Item("Главные")
$wsLinear = $wb.Worksheets.Item("Линейные")

# --- Главные: updated statistic values ---
$wsMainUpdates = @(
    @{ Row=3; Cells=@{ "C"="27"; "D"="480"; "E"="221"; "F"="259"; "G"="17.78"; "H"="8.19"; "I"="9.59"; "J"="108"; "K"="107" } }
    @{ Row=4; Cells=@{ "C"="21"; "D"="334"; "E"="143"; "F"="191"; "G"="15.9"; "H"="6.81"; "I"="9.1"; "J"="69"; "K"="83"; "V"="2" } }
    @{ Row=5; Cells=@{ "C"="27"; "D"="433"; "E"="226"; "F"="207"; "G"="16.04"; "H"="8.369999999999999"; "I"="7.67"; "J"="108"; "K"="96"; "V"="22"; "W"="18" } }
    @{ Row=6; Cells=@{ "C"="27"; "D"="457"; "E"="201"; "F"="256"; "G"="16.93"; "H"="7.44"; "I"="9.48"; "J"="93"; "K"="108" } }
    @{ Row=8; Cells=@{ "C"="24"; "D"="434"; "E"="219"; "F"="215"; "G"="18.08"; "H"="9.130000000000001"; "I"="8.960000000000001"; "J"="102"; "K"="100"; "V"="12"; "W"="12" } }
    @{ Row=12; Cells=@{ "C"="18"; "D"="307"; "E"="133"; "F"="174"; "G"="17.06"; "H"="7.39"; "I"="9.67"; "J"="54"; "K"="62"; "V"="10"; "W"="6" } }
    @{ Row=19; Cells=@{ "C"="22"; "D"="386"; "E"="194"; "F"="192"; "G"="17.55"; "H"="8.82"; "I"="8.73"; "J"="92"; "K"="81" } }
    @{ Row=20; Cells=@{ "C"="25"; "D"="422"; "E"="180"; "F"="242"; "G"="16.88"; "H"="7.2"; "I"="9.68"; "J"="85"; "K"="91" } }
    @{ Row=23; Cells=@{ "C"="17"; "D"="218"; "E"="81"; "F"="137"; "G"="12.82"; "H"="4.76"; "I"="8.06"; "J"="38"; "K"="56"; "V"="4" } }
    @{ Row=24; Cells=@{ "C"="28"; "D"="463"; "E"="209"; "F"="254"; "G"="16.54"; "H"="7.46"; "I"="9.07"; "J"="102"; "K"="112"; "V"="14"; "W"="12" } }
)

foreach ($update in $wsMainUpdates) {
    $r = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $wsMain.Range("$col$r").Value = [double]$update.Cells[$col]
    }
}

# --- Линейные: updated statistic values ---
$wsLinearUpdates = @(
    @{ Row=3; Cells=@{ "C"="26"; "D"="381"; "E"="197"; "F"="184"; "G"="14.65"; "H"="7.58"; "I"="7.08"; "J"="96"; "K"="77" } }
    @{ Row=5; Cells=@{ "C"="14"; "D"="202"; "E"="108"; "F"="94"; "G"="14.43"; "H"="7.71"; "I"="6.71"; "J"="54"; "K"="47"; "V"="12"; "W"="6" } }
    @{ Row=7; Cells=@{ "C"="17"; "D"="271"; "E"="96"; "F"="175"; "G"="15.94"; "H"="5.65"; "I"="10.29"; "J"="48"; "K"="60"; "V"="10"; "W"="8" } }
    @{ Row=8; Cells=@{ "C"="25"; "D"="387"; "E"="148"; "F"="239"; "G"="15.48"; "H"="5.92"; "I"="9.56"; "J"="69"; "K"="92"; "V"="6" } }
    @{ Row=9; Cells=@{ "C"="25"; "D"="462"; "E"="201"; "F"="261"; "G"="18.48"; "H"="8.039999999999999"; "I"="10.44"; "J"="88"; "K"="108"; "V"="12" } }
    @{ Row=14; Cells=@{ "C"="25"; "D"="412"; "E"="210"; "F"="202"; "G"="16.48"; "H"="8.4"; "I"="8.08"; "J"="105"; "K"="96" } }
    @{ Row=15; Cells=@{ "C"="23"; "D"="435"; "E"="227"; "F"="208"; "G"="18.91"; "H"="9.869999999999999"; "I"="9.039999999999999"; "J"="91"; "K"="84"; "V"="12"; "W"="12" } }
    @{ Row=19; Cells=@{ "C"="24"; "D"="415"; "E"="200"; "F"="215"; "G"="17.29"; "H"="8.33"; "I"="8.960000000000001"; "J"="95"; "K"="95" } }
    @{ Row=20; Cells=@{ "C"="18"; "D"="288"; "E"="147"; "F"="141"; "G"="16"; "H"="8.17"; "I"="7.83"; "J"="71"; "K"="68" } }
)

foreach ($update in $wsLinearUpdates) {
    $r = $update.Row
    foreach ($col in $update.Cells.Keys) {
        $wsLinear.Range("$col$r").Value = [double]$update.Cells[$col]
    }
}

# --- Главные: refresh AA timestamp for every data row (2-26) ---
for ($r = 2; $r -le 26; $r++) {
    $wsMain.Range("AA$r").Value = $timestamp
}

# --- Линейные: refresh AA timestamp for every data row (2-26) ---
for ($r = 2; $r -le 26; $r++) {
    $wsLinear.Range("AA$r").Value = $timestamp
}
